$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 273, shifting existing rows 273-381 down to 274-382
$ws.Rows.Item(273).Insert()

# Populate the new row 273 with the new data
$ws.Cells.Item(273, 1).Value = 5
$ws.Cells.Item(273, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(273, 3).Value = "Maule"
$ws.Cells.Item(273, 4).Value = 44755
$ws.Cells.Item(273, 5).Value = 7
$ws.Cells.Item(273, 6).Value = 100112032
$ws.Cells.Item(273, 7).Value = "Zapallo italiano"
$ws.Cells.Item(273, 8).Value = "Sin especificar"
$ws.Cells.Item(273, 9).Value = "Primera"
$ws.Cells.Item(273, 10).Value = 300
$ws.Cells.Item(273, 11).Value = 12000
$ws.Cells.Item(273, 12).Value = 12000
$ws.Cells.Item(273, 13).Value = 12000
$ws.Cells.Item(273, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(273, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(273, 16).Value = 240
$ws.Cells.Item(273, 17).Value = 50
$ws.Cells.Item(273, 18).Value = "Hortaliza"

# Apply the date style to column D (matches the other rows' date-number-format cell style)
$ws.Cells.Item(273, 4).NumberFormat = $ws.Cells.Item(274, 4).NumberFormat
